$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3218.182
$ws.Range("I76").Value = 3214.2856
$ws.Range("J76").Value = 3225
$ws.Range("K76").Value = 3214.2856
$ws.Range("L76").Value = 3225
$ws.Range("M76").Value = -2899.2856
$ws.Range("N76").Value = -3855

$ws.Range("H79").Value = 3218.182
$ws.Range("I79").Value = 3214.2856
$ws.Range("J79").Value = 3225
$ws.Range("K79").Value = 3214.2856
$ws.Range("L79").Value = 3225
$ws.Range("M79").Value = -2122.2856
$ws.Range("N79").Value = -5409

$ws.Range("H100").Value = 1128.0714
$ws.Range("I100").Value = 1183.1666
$ws.Range("J100").Value = 797.5
$ws.Range("K100").Value = 1183.1666
$ws.Range("L100").Value = 797.5
$ws.Range("M100").Value = -642.1666
$ws.Range("N100").Value = -1879.5

$ws.Range("H116").Value = 1940.5
$ws.Range("I116").Value = 1401
$ws.Range("J116").Value = 2480
$ws.Range("K116").Value = 1401
$ws.Range("L116").Value = 2480
$ws.Range("M116").Value = 2041
$ws.Range("N116").Value = -9364

$ws.Range("H129").Value = 938.9286
$ws.Range("I129").Value = 568.3333
$ws.Range("J129").Value = 1000.69446
$ws.Range("K129").Value = 1704.9999
$ws.Range("L129").Value = 3002.08338
$ws.Range("M129").Value = 3295.0001
$ws.Range("N129").Value = -13002.08338

$ws.Range("H141").Value = 1501.4546
$ws.Range("I141").Value = 843.7895
$ws.Range("J141").Value = 5666.6665
$ws.Range("K141").Value = 2531.3685
$ws.Range("L141").Value = 16999.9995
$ws.Range("M141").Value = 2648.6315
$ws.Range("N141").Value = -27359.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24158.785
$ws.Range("I32").Value = 26349.695
$ws.Range("K32").Value = 26349.695
$ws.Range("M32").Value = -26062.695

$ws.Range("H110").Value = 1696.2858
$ws.Range("I110").Value = 1719.5
$ws.Range("J110").Value = 1665.3334
$ws.Range("K110").Value = 1719.5
$ws.Range("L110").Value = 1665.3334
$ws.Range("M110").Value = 325.5
$ws.Range("N110").Value = -5755.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 41668176
$ws.Range("I105").Value = 62501200
$ws.Range("J105").Value = 2125
$ws.Range("K105").Value = 62501200
$ws.Range("L105").Value = 2125
$ws.Range("M105").Value = -62499453
$ws.Range("N105").Value = -5619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2162.6
$ws.Range("I16").Value = 1600
$ws.Range("J16").Value = 2303.25
$ws.Range("K16").Value = 1600
$ws.Range("L16").Value = 2303.25
$ws.Range("M16").Value = -1313
$ws.Range("N16").Value = -2877.25

$ws.Range("H31").Value = 2718.4055
$ws.Range("I31").Value = 1354.5518
$ws.Range("J31").Value = 7662.375
$ws.Range("K31").Value = 1354.5518
$ws.Range("L31").Value = 7662.375
$ws.Range("M31").Value = -1059.5518
$ws.Range("N31").Value = -8252.375

$ws.Range("H34").Value = 2718.4055
$ws.Range("I34").Value = 1354.5518
$ws.Range("J34").Value = 7662.375
$ws.Range("K34").Value = 1354.5518
$ws.Range("L34").Value = 7662.375
$ws.Range("M34").Value = -1152.5518
$ws.Range("N34").Value = -8066.375

$ws.Range("H105").Value = 786.7
$ws.Range("I105").Value = 764
$ws.Range("J105").Value = 877.5
$ws.Range("K105").Value = 764
$ws.Range("L105").Value = 877.5
$ws.Range("M105").Value = 983
$ws.Range("N105").Value = -4371.5

$ws.Range("H107").Value = 427.8095
$ws.Range("I107").Value = 369.64706
$ws.Range("J107").Value = 675
$ws.Range("K107").Value = 369.64706
$ws.Range("L107").Value = 675
$ws.Range("M107").Value = 1550.35294
$ws.Range("N107").Value = -4515

$ws.Range("H113").Value = 2162.6
$ws.Range("I113").Value = 1600
$ws.Range("J113").Value = 2303.25
$ws.Range("K113").Value = 1600
$ws.Range("L113").Value = 2303.25
$ws.Range("M113").Value = 570
$ws.Range("N113").Value = -6643.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.14286
$ws.Range("I2").Value = 104.2
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 625.2
$ws.Range("L2").Value = 540
$ws.Range("M2").Value = -512.2
$ws.Range("N2").Value = -766

$ws.Range("H22").Value = 878.625
$ws.Range("I22").Value = 590
$ws.Range("J22").Value = 974.8333
$ws.Range("K22").Value = 1770
$ws.Range("L22").Value = 2924.4999
$ws.Range("M22").Value = -1601
$ws.Range("N22").Value = -3262.4999

$ws.Range("H27").Value = 878.625
$ws.Range("I27").Value = 590
$ws.Range("J27").Value = 974.8333
$ws.Range("K27").Value = 1770
$ws.Range("L27").Value = 2924.4999
$ws.Range("M27").Value = -1668
$ws.Range("N27").Value = -3128.4999

$ws.Range("H113").Value = 562.7162
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 622.1316
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 1866.3948
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6206.3948

$ws.Range("H121").Value = 45445390
$ws.Range("J121").Value = 53464944
$ws.Range("L121").Value = 160394832
$ws.Range("N121").Value = -160397452

$ws.Range("H122").Value = 841.8570999999999
$ws.Range("I122").Value = 286.75
$ws.Range("K122").Value = 2580.75
$ws.Range("M122").Value = -130.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27584.838
$ws.Range("I70").Value = 36201.547
$ws.Range("J70").Value = 5325
$ws.Range("K70").Value = 36201.547
$ws.Range("L70").Value = 5325
$ws.Range("M70").Value = -35931.547
$ws.Range("N70").Value = -5865

$ws.Range("H73").Value = 27584.838
$ws.Range("I73").Value = 36201.547
$ws.Range("J73").Value = 5325
$ws.Range("K73").Value = 36201.547
$ws.Range("L73").Value = 5325
$ws.Range("M73").Value = -35265.547
$ws.Range("N73").Value = -7197

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 38130.4
$ws.Range("J70").Value = 38130.4
$ws.Range("L70").Value = 38130.4
$ws.Range("N70").Value = -38670.4

$ws.Range("H73").Value = 38130.4
$ws.Range("J73").Value = 38130.4
$ws.Range("L73").Value = 38130.4
$ws.Range("N73").Value = -40002.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 21500
$ws.Range("J63").Value = 21500
$ws.Range("L63").Value = 21500
$ws.Range("N63").Value = -22748

$ws.Range("H66").Value = 21500
$ws.Range("J66").Value = 21500
$ws.Range("L66").Value = 64500
$ws.Range("N66").Value = -70740

$ws.Range("H100").Value = 101169.6
$ws.Range("I100").Value = 250550
$ws.Range("J100").Value = 63824.5
$ws.Range("K100").Value = 501100
$ws.Range("L100").Value = 127649
$ws.Range("M100").Value = -500559
$ws.Range("N100").Value = -128731

Write-Output "edits applied"
